$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

# "Recorded By" column (G) - reorder the author list on every row that
# currently reads "<email>, System" to "System, <email>".
$colG = $ws.Range("G1:G259")

$first = $colG.Find($oldText)
if ($first -ne $null) {
    $firstAddr = $first.Address()
    $current = $first
    do {
        $current.Value = $newText
        $current = $colG.FindNext($current)
    } while ($current -ne $null -and $current.Address() -ne $firstAddr)
}
